$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Clear the "Liens Ressources Humaines :" / "K'IAM" block that used to repeat
# in rows 10, 13, 25, 31 and 36 (A:D), leaving only the formatting behind.
$ws.Range("A10:D10").ClearContents()
$ws.Range("A13:D13").ClearContents()
$ws.Range("A25:D25").ClearContents()
$ws.Range("A31:D31").ClearContents()
$ws.Range("A36:D36").ClearContents()

# Update the active selection to reflect the last edited cell.
$ws.Range("D36").Select()
